# Apply the "Cập nhật và sửa lỗi cac chuc nang moi nhat" update to the
# evaluation summary worksheet.
#
# Net content changes on the single data worksheet:
#   - E5 gets a new numeric value (62)
#   - M8 changes from "A" to "C"
#   - K10 changes from 50 to 40
#   - K11 changes from 85 to 100
#   - I12 changes from "A" to "C"
#   - J12 changes from 74 to 97.5
#   - K12 changes from 79.5 to 85.2
#   - L12 changes from 100 to 75
#   - M12 changes from "A" to "C"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: fill in the previously blank "Số ngày làm việc thực tế" value.
$ws.Range("E5").Value = 62

# Row 8: "Mức xếp loại của Lãnh đạo" changes from A to C.
$ws.Range("M8").Value = "C"

# Row 10: updated score.
$ws.Range("K10").Value = 40

# Row 11: updated score.
$ws.Range("K11").Value = 100

# Row 12: rating changes from A to C; scores updated.
$ws.Range("I12").Value = "C"
$ws.Range("J12").Value = 97.5
$ws.Range("K12").Value = 85.2
$ws.Range("L12").Value = 75
$ws.Range("M12").Value = "C"

$wb.Save()
